$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header fields
$ws.Range("G4").Value = 44079.51026330685
$ws.Range("G7").Value = "Unesco"
$ws.Range("F10").Value = "test"

# Row 18 - piqlConnect (only piqlFilm)
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1500
$ws.Range("H18").Value = 1500

# Row 19 - Digital (GB) - clear
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""

# Row 20 - Visual (pages) - E20 must be stored as text "1" (not a number),
# while keeping the cell's existing numeric style/format (0.00). Temporarily
# switch to a text format so Excel doesn't coerce "1" into a number, then
# restore the original format.
$origFormat = $ws.Range("E20").NumberFormat
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1"
$ws.Range("E20").NumberFormat = $origFormat
$ws.Range("F20").Value = 67000
$ws.Range("G20").Value = 0.03
$ws.Range("H20").Value = 2010

# Row 21 - Online Storage (GB) - piqlConnect (1TB & 1 piqlFilm included) - clear
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = ""

# Row 22 - Online Storage (GB), Payment - clear
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = ""
$ws.Range("G22").Value = ""
$ws.Range("H22").Value = ""

# Row 24 - Registration fee - clear
$ws.Range("F24").Value = ""
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = ""

# Row 25 - AWA contribution, Entity - clear
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = ""

# Row 26 - Management fee (per year) - clear
$ws.Range("F26").Value = ""
$ws.Range("G26").Value = ""
$ws.Range("H26").Value = ""

# Row 27 - Storage (reels/ per year), Period (years) - clear
$ws.Range("E27").Value = ""
$ws.Range("F27").Value = ""
$ws.Range("G27").Value = ""
$ws.Range("H27").Value = ""

# Row 29 - piqlReader - clear
$ws.Range("F29").Value = ""
$ws.Range("G29").Value = ""
$ws.Range("H29").Value = ""

# Row 30 - Installation and training - clear
$ws.Range("F30").Value = ""
$ws.Range("G30").Value = ""
$ws.Range("H30").Value = ""

# Row 31 - Service agreement (per year), Type - clear
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = ""
$ws.Range("G31").Value = ""
$ws.Range("H31").Value = ""

# Row 32 - Shipment cost, Reels
$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 30
$ws.Range("H32").Value = 60

# Row 33 - TOTAL
$ws.Range("H33").Value = 3510

# Row 34 - Total to pay from the second term
$ws.Range("H34").Value = 0
